# Fixed update to excel issue
# - Rename "Requested quantity" header -> "Weekly_PO_Qty" on "Weekly Quantity" sheet
# - Rename "Requested quantity" header -> "Monthly_PO_Qty" on "Monthly Trend" sheet
# - Add a new "PO Forecast" sheet with ds / PO_Forecast / yhat_lower / yhat_upper data

$wb = $excel.ActiveWorkbook

# --- Rename headers on the existing sheets -------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet -------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Re-use the existing header formatting (bold, centered, top-aligned, thin
# border) by copying the format from the "Weekly Quantity" header row, then
# overwrite the text.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsForecast.Cells.Item(1, 1).Value = "ds"
$wsForecast.Cells.Item(1, 2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1, 3).Value = "yhat_lower"
$wsForecast.Cells.Item(1, 4).Value = "yhat_upper"

# Re-use the existing date-column formatting for column A's data rows.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A18").PasteSpecial(-4122)

# Data rows: ds, PO_Forecast, yhat_lower, yhat_upper
$data = @(
    @(45431.99999999999, 6, 1.705939123019065, 10.17312013948997),
    @(45438.99999999999, 6, 1.293833826129422, 10.25007867866567),
    @(45445.99999999999, 6, 1.0783435354472, 10.40868489082245),
    @(45452.99999999999, 5, 0.7023779044422994, 9.978603841433859),
    @(45473.99999999999, 5, 0.7759659883690634, 9.344935521835357),
    @(45557.99999999999, 3, -1.519859954419893, 7.540813136339012),
    @(45592.99999999999, 2, -2.73504594001737, 6.790436688585587),
    @(45599.99999999999, 2, -2.575867173673854, 6.282706720807741),
    @(45606.99999999999, 2, -2.686668900611354, 6.284309357678722),
    @(45613.99999999999, 1, -2.8881585763257, 6.290980036581488),
    @(45620.99999999999, 1, -3.237939200132748, 5.773665503950388),
    @(45627.99999999999, 1, -3.14027711161347, 5.882475059971522),
    @(45634.99999999999, 1, -3.541647691175124, 5.636058348364071),
    @(45641.99999999999, 1, -3.56948908107225, 5.305857415789937),
    @(45648.99999999999, 1, -4.041764262725526, 4.877974270909896),
    @(45655.99999999999, 0, -3.640110812626593, 4.89375541795898),
    @(45662.99999999999, 0, -3.934445247411893, 4.887580110968679)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $row = $data[$i]
    $wsForecast.Cells.Item($rowNum, 1).Value = $row[0]
    $wsForecast.Cells.Item($rowNum, 2).Value = $row[1]
    $wsForecast.Cells.Item($rowNum, 3).Value = $row[2]
    $wsForecast.Cells.Item($rowNum, 4).Value = $row[3]
}

Write-Output "PO Forecast sheet added; headers renamed."
